# "create criteria and tension tests and implementation"
#
# 1. Rename Sheet1 -> "W profile slenderness" and move its selection to E4.
# 2. Insert a new "Tension" sheet right after it, with headers + the
#    tension-capacity calculation (gross/effective/nominal/LRFD values for
#    a w6x15 profile), selection left on M3 as the active sheet/cell.

$wb = $excel.ActiveWorkbook

# --- Rename & reselect the existing sheet -----------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "W profile slenderness"
$ws1.Range("E4").Select() | Out-Null

# --- Add the new "Tension" sheet right after it ------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Tension"

# Headers (row 1)
$ws2.Range("A1").Value = "yield strength Mpa"
$ws2.Range("B1").Value = "ultimate strength Mpa"
$ws2.Range("C1").Value = "Profile"
$ws2.Range("D1").Value = "gross area"
$ws2.Range("E1").Value = "net area"
$ws2.Range("F1").Value = "lag factor "
$ws2.Range("G1").Value = "effective are"
$ws2.Range("H1").Value = "nominal ult str"
$ws2.Range("I1").Value = "nominal yield strength"
$ws2.Range("J1").Value = "ultimate strength lrfd"
$ws2.Range("K1").Value = "yield strength lrfd"

# Data (row 2)
$ws2.Range("A2").Value = 355
$ws2.Range("B2").Value = 500
$ws2.Range("C2").Value = "w6x15"
$ws2.Range("D2").Value = 2860
$ws2.Range("E2").Value = 2860
$ws2.Range("F2").Value = 1
$ws2.Range("G2").Formula = "=F2*E2"
$ws2.Range("H2").Formula = "=G2*B2"
$ws2.Range("I2").Formula = "=D2*A2"
$ws2.Range("J2").Formula = "=H2*0.75"
$ws2.Range("K2").Formula = "=0.9*I2"
$ws2.Range("L2").Formula = "=H2/2"
$ws2.Range("M2").Formula = "=I2/1.67"

# Leave the Tension sheet active, selection on M3 (matches the target state)
$ws2.Range("M3").Select() | Out-Null
